$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cargos")
Write-Output $ws.Name
